$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the whole used range first so stale columns C/D and row 7 are removed
$ws.Cells.Clear()

# New header row
$ws.Range("A1").Value = "Casos de uso "

# Shifted data rows (now only columns A and B are used)
$ws.Range("A2").Value = "Crear Solicitud"
$ws.Range("B2").Value = "Pag 44"

$ws.Range("A3").Value = "Agregar Cotizacion"
$ws.Range("B3").Value = "Pag 54"

$ws.Range("A4").Value = "Solicitar Partida Especial"
$ws.Range("B4").Value = "Pag 59"

$ws.Range("A5").Value = "Registrar Adquisición"
$ws.Range("B5").Value = "Pag 81"

$ws.Range("A6").Value = "Crear Asignacion"
$ws.Range("B6").Value = "Pag 117"

# Row 7 intentionally left blank

$ws.Range("A8").Value = "Casos de prueba"
$ws.Range("B8").Value = "Pag 180"

# Update selection to match target state
$ws.Range("A10").Select() | Out-Null

$wb.Save()
